$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.388442635536194
$ws.Range("B1").Value = 3.437412261962891
$ws.Range("C1").Value = 3.438670635223389
$ws.Range("D1").Value = 1.620058655738831
$ws.Range("E1").Value = 1.227767705917358
